$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# Row 11 - TP7 ("--", nota 0)
# ---------------------------------------------------------------
$ws.Range("A7").Copy()
$ws.Range("A11").PasteSpecial(-4122)
$ws.Range("A11").Value = "TP7"

$ws.Range("B7").Copy()
$ws.Range("B11").PasteSpecial(-4122)
$ws.Range("B11").Value = "'--"

$ws.Range("C7").Copy()
$ws.Range("C11").PasteSpecial(-4122)
$ws.Range("C11").Value = 0

$ws.Rows.Item(11).RowHeight = 17

# ---------------------------------------------------------------
# Row 12 - TP8 (Clusterização..., nota (0+6)/2, *Reavaliado)
# ---------------------------------------------------------------
$ws.Range("A8").Copy()
$ws.Range("A12").PasteSpecial(-4122)
$ws.Range("A12").Value = "TP8"

$ws.Range("B8").Copy()
$ws.Range("B12").PasteSpecial(-4122)
$ws.Range("B12").Value = "Clusterização em projetos Spring Boot`n- anlr4 a partir de python`n- sairam do Python e foram para o java.`n- KNN modificado (qual a modificação?)`n- Não entendi bem o que faz conceitualmente. Isso tem que melhorar nos próximos TPs.`n- clusterização a partir de dados principais (importações, etc.)`n- Perdeu muito tempo com explicação do código, mas faltou entender o que de fato está fazendo. No final falou repository, service, controller, model...`n- Limitação forte do @Entity, mas entendo."

$ws.Range("C8").Copy()
$ws.Range("C12").PasteSpecial(-4122)
$ws.Range("C12").Formula = "=(0+6)/2"

$ws.Range("D6").Copy()
$ws.Range("D12").PasteSpecial(-4122)
$ws.Range("D12").Value = "*Reavaliado"

$ws.Rows.Item(12).RowHeight = 204

# ---------------------------------------------------------------
# Row 13 - TP9 (Descobridor de arquitetura MVC e MVP..., nota 4)
# ---------------------------------------------------------------
$ws.Range("A7").Copy()
$ws.Range("A13").PasteSpecial(-4122)
$ws.Range("A13").Value = "TP9"

$ws.Range("B11").Copy()
$ws.Range("B13").PasteSpecial(-4122)
$ws.Range("B13").Value = "'Descobridor de arquitetura MVC e MVP`n- está mais para padrão arquitetural do que para descoberta arquitetural.`n- Antlr4 + java`n- Ficou muito aquém do que se espera de `"descoberta arquitetural`". Na verdade, muito próximo ao TP8 de vocês focado em padrão arquitetural. Obsersem que descoberta arquitetural, você quer ter uma visão global de como a aplicação funciona e não apenas enquadrar em um padrão arquitetural.`n- Lembrar de não focar tanto no código, mas na explicação teórica do que você faz. Nesse eu entendi, mas no TP8 ficou complicado.`n- Bom apontar limitações e dificuldades."

$ws.Range("C7").Copy()
$ws.Range("C13").PasteSpecial(-4122)
$ws.Range("C13").Value = 4

$ws.Rows.Item(13).RowHeight = 238

# ---------------------------------------------------------------
# Selection
# ---------------------------------------------------------------
$ws.Range("B14").Select()
